$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 7484.2354
$ws.Range("I40").Value = 5833.7
$ws.Range("K40").Value = 5833.7
$ws.Range("M40").Value = -5658.7
$ws.Range("H51").Value = 2500
$ws.Range("I51").Value = 2500
$ws.Range("K51").Value = 2500
$ws.Range("M51").Value = -2016
$ws.Range("H69").Value = 7195.8066
$ws.Range("J69").Value = 7195.8066
$ws.Range("L69").Value = 21587.4198
$ws.Range("N69").Value = -23335.4198
$ws.Range("H72").Value = 7195.8066
$ws.Range("J72").Value = 7195.8066
$ws.Range("L72").Value = 64762.2594
$ws.Range("N72").Value = -73498.2594
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H86").Value = 6585.875
$ws.Range("I86").Value = 6615.6665
$ws.Range("J86").Value = 6496.5
$ws.Range("K86").Value = 6615.6665
$ws.Range("L86").Value = 6496.5
$ws.Range("M86").Value = -5492.6665
$ws.Range("N86").Value = -8742.5
$ws.Range("H89").Value = 6585.875
$ws.Range("I89").Value = 6615.6665
$ws.Range("J89").Value = 6496.5
$ws.Range("K89").Value = 33078.3325
$ws.Range("L89").Value = 32482.5
$ws.Range("M89").Value = -27462.3325
$ws.Range("N89").Value = -43714.5
$ws.Range("H92").Value = 145.66667
$ws.Range("I92").Value = 145.66667
$ws.Range("K92").Value = 145.66667
$ws.Range("M92").Value = 1102.33333
$ws.Range("H105").Value = 45000
$ws.Range("J105").Value = 45000
$ws.Range("L105").Value = 45000
$ws.Range("N105").Value = -51988

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3955.5
$ws.Range("J61").Value = 8999.75
$ws.Range("L61").Value = 8999.75
$ws.Range("N61").Value = -9423.75
$ws.Range("H74").Value = 5573.615
$ws.Range("I74").Value = 5563.778
$ws.Range("J74").Value = 5595.75
$ws.Range("K74").Value = 5563.778
$ws.Range("L74").Value = 5595.75
$ws.Range("M74").Value = -4689.778
$ws.Range("N74").Value = -7343.75
$ws.Range("H77").Value = 5573.615
$ws.Range("I77").Value = 5563.778
$ws.Range("J77").Value = 5595.75
$ws.Range("K77").Value = 27818.89
$ws.Range("L77").Value = 27978.75
$ws.Range("M77").Value = -23450.89
$ws.Range("N77").Value = -36714.75
$ws.Range("H88").Value = 983
$ws.Range("J88").Value = 634.3333
$ws.Range("L88").Value = 634.3333
$ws.Range("N88").Value = -1446.3333
$ws.Range("H91").Value = 983
$ws.Range("J91").Value = 634.3333
$ws.Range("L91").Value = 634.3333
$ws.Range("N91").Value = -3442.3333
$ws.Range("H92").Value = 25000
$ws.Range("J92").Value = 25000
$ws.Range("L92").Value = 25000
$ws.Range("N92").Value = -29992
$ws.Range("H132").Value = 5316.25
$ws.Range("I132").Value = 5316.25
$ws.Range("K132").Value = 15948.75
$ws.Range("M132").Value = -13418.75
$ws.Range("H136").Value = 3955.5
$ws.Range("J136").Value = 8999.75
$ws.Range("L136").Value = 26999.25
$ws.Range("N136").Value = -32099.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1911.7858
$ws.Range("I22").Value = 1911.7858
$ws.Range("K22").Value = 1911.7858
$ws.Range("M22").Value = -1738.7858
$ws.Range("H86").Value = 3648.7827
$ws.Range("I86").Value = 1813.5
$ws.Range("J86").Value = 5650.909
$ws.Range("K86").Value = 1813.5
$ws.Range("L86").Value = 5650.909
$ws.Range("M86").Value = -690.5
$ws.Range("N86").Value = -7896.909
$ws.Range("H89").Value = 3648.7827
$ws.Range("I89").Value = 1813.5
$ws.Range("J89").Value = 5650.909
$ws.Range("K89").Value = 9067.5
$ws.Range("L89").Value = 28254.545
$ws.Range("M89").Value = -3451.5
$ws.Range("N89").Value = -39486.545

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 3450.0322
$ws.Range("I7").Value = 3925.1155
$ws.Range("K7").Value = 3925.1155
$ws.Range("M7").Value = -3812.1155
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()
$ws.Range("H92").Value = 30861.6
$ws.Range("J92").Value = 30861.6
$ws.Range("L92").Value = 30861.6
$ws.Range("N92").Value = -35853.6
$ws.Range("H98").Value = 79890
$ws.Range("J98").Value = 79890
$ws.Range("L98").Value = 79890
$ws.Range("N98").Value = -84382

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("M92").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 10067891
$ws.Range("J11").Value = 5557833
$ws.Range("L11").Value = 5557833
$ws.Range("N11").Value = -5558111

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H26").Value = 1833
$ws.Range("I26").Value = 2249.5
$ws.Range("K26").Value = 2249.5
$ws.Range("M26").Value = -1954.5
$ws.Range("H55").Value = 845.7647
$ws.Range("I55").Value = 758.5333
$ws.Range("J55").Value = 1500
$ws.Range("K55").Value = 758.5333
$ws.Range("L55").Value = 1500
$ws.Range("M55").Value = -585.5333
$ws.Range("N55").Value = -1846
$ws.Range("H68").Value = 7084.4
$ws.Range("J68").Value = 10000
$ws.Range("L68").Value = 10000
$ws.Range("N68").Value = -11498
$ws.Range("H71").Value = 7084.4
$ws.Range("J71").Value = 10000
$ws.Range("L71").Value = 50000
$ws.Range("N71").Value = -57488
$ws.Range("H100").Value = 6541
$ws.Range("I100").Value = 3098.6
$ws.Range("J100").Value = 8999.857
$ws.Range("K100").Value = 3098.6
$ws.Range("L100").Value = 8999.857
$ws.Range("M100").Value = -2557.6
$ws.Range("N100").Value = -10081.857
$ws.Range("H122").Value = 2550.2942
$ws.Range("I122").Value = 2167.125
$ws.Range("J122").Value = 2890.889
$ws.Range("K122").Value = 6501.375
$ws.Range("L122").Value = 8672.667000000001
$ws.Range("M122").Value = -4051.375
$ws.Range("N122").Value = -13572.667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5799.6
$ws.Range("I81").Value = 2999.6667
$ws.Range("J81").Value = 9999.5
$ws.Range("K81").Value = 5999.3334
$ws.Range("L81").Value = 19999
$ws.Range("M81").Value = -4938.3334
$ws.Range("N81").Value = -22121
$ws.Range("H84").Value = 5799.6
$ws.Range("I84").Value = 2999.6667
$ws.Range("J84").Value = 9999.5
$ws.Range("K84").Value = 29996.667
$ws.Range("L84").Value = 99995
$ws.Range("M84").Value = -24692.667
$ws.Range("N84").Value = -110603
